$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "name" field right after "campaign_id" (row 3).
# This pushes total_redemption_count / total_redemption_amount /
# active_coupon_code_count / created_at / updated_at down by one row, while
# keeping their formatting/content intact (no further edits needed for the
# rows that land on 5, 6 and 7).
$ws.Rows.Item(4).Insert()

# The insert copied the "UUID" example cell from J3 down into J4; the new
# "name" row doesn't use column J, so drop that stray cell entirely.
$ws.Range("J4").Clear()

# Row 4: name (use A3's normal style so every cell matches the sheet's
# regular formatting instead of inheriting K3's shaded "examples" style).
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A4").Value = "name"
$ws.Range("A3").Copy($ws.Range("B4"))
$ws.Range("B4").Value = "No"
$ws.Range("A3").Copy($ws.Range("C4"))
$ws.Range("C4").Value = "string"
$ws.Range("A3").Copy($ws.Range("K4"))
$ws.Range("K4").Value = "default"

# Rows 5-7 (total_redemption_count, total_redemption_amount,
# active_coupon_code_count) keep their previous content unchanged - nothing
# else to do for them.

# Row 8 used to hold "created_at"; it now becomes the new "status" field.
$ws.Range("A8").Value = "status"
$ws.Range("B8").Value = "No"
$ws.Range("C8").Value = "string"
$ws.Range("A3").Copy($ws.Range("I8"))
$ws.Range("I8").Value = "active, inactive"
$ws.Range("K8").Value = "active"

# Row 9 used to hold "updated_at"; it now becomes "created_at".
$ws.Range("A9").Value = "created_at"
$ws.Range("B9").Value = "No"
$ws.Range("C9").Value = "string"
$ws.Range("K9").Value = "2025-03-12 20:24:03.653247+05:30"

# Row 10 is a brand new row for "updated_at"; copy styles from row 9 so its
# cells carry the sheet's normal formatting.
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "updated_at"
$ws.Range("B9").Copy($ws.Range("B10"))
$ws.Range("B10").Value = "No"
$ws.Range("C9").Copy($ws.Range("C10"))
$ws.Range("C10").Value = "string"
$ws.Range("K9").Copy($ws.Range("K10"))
$ws.Range("K10").Value = "2025-03-12 20:24:03.653247+05:30"
